# Updated data for aged orders
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Data")

# Add new aged order row 8
$ws.Range("B8").Value = 60588621
$ws.Range("F8").Value = "STARTECH.COM"

# Update BCN / order number for existing aged order row (row 2)
$ws.Range("B2").Value = 30036215

# Update the vendor/product name for row 2
$ws.Range("F2").Value = "APPLE MINI IPADS"

# Match the style used on B2 for the newly added B8 cell
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Reflect the new active cell selection
$ws.Range("I13").Select()
